$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Ajo" (Terminal
# Hortofrutícola Agro Chillán). Insert a new row at position 68, which
# pushes the existing rows 68-137 down to 69-138, and fill the new row 68
# with the new observation's data. Columns A,B,C,E,F,G,H,I,N,O,Q,R keep the
# same values as the (old) row 68 / the rest of the table.

$ws.Rows(68).Insert()

$ws.Range("A68").Value = 7
$ws.Range("B68").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C68").Value = "Ñuble"
$ws.Range("D68").Value = 44494
$ws.Range("E68").Value = 16
$ws.Range("F68").Value = 100112003
$ws.Range("G68").Value = "Ajo"
$ws.Range("H68").Value = "Chino"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 100
$ws.Range("K68").Value = 16000
$ws.Range("L68").Value = 17000
$ws.Range("M68").Value = 16500
$ws.Range("N68").Value = "$/caja 10 kilos"
$ws.Range("O68").Value = "China"
$ws.Range("P68").Value = 1650
$ws.Range("Q68").Value = 10
$ws.Range("R68").Value = "Hortaliza"
